$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "object_detection_zed_" + "live" -> single run "object_detection_zed_live"
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Find.Execute("object_detection_zed_live", $false, $false, $false, $false, $false, $true, 1, $false, "object_detection_zed_live", 2)

# ---------------------------------------------------------------------
# 2) "object_detection_zed_" + "modi" -> single run "object_detection_zed_modi"
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Find.Execute("object_detection_zed_modi", $false, $false, $false, $false, $false, $true, 1, $false, "object_detection_zed_modi", 2)

# ---------------------------------------------------------------------
# 3) "Saving the output in csv with real world coordinates" gets split
#    into "Saving the ou" | _GoBack bookmark | "tput in csv with real
#    world coordinates" (the _GoBack bookmark simply marks the last
#    edited spot, it moves here from the end of the document).
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$rngSplit = $p4.Range
$rngSplit.Find.Execute("Saving the ou", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngSplit.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rngSplit) | Out-Null

# ---------------------------------------------------------------------
# 4) "Results from " + "print(type(" + "x" + "_" + <bookmark _Hlk16155600> "center" </bookmark>
#    -> "Results from print(type(" + "x_" + <bookmark> "center" </bookmark>
# ---------------------------------------------------------------------
$p22 = $d.Paragraphs.Item(22)
$p22.Range.Find.Execute("Results from print(type(", $false, $false, $false, $false, $false, $true, 1, $false, "Results from print(type(", 2)

# ---------------------------------------------------------------------
# 5) second "x" + "_center" -> "x_center" (the one right after "print(")
# ---------------------------------------------------------------------
$markRng = $p22.Range
$markRng.Find.Execute(")), print(", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$afterPrint = $d.Range($markRng.End, $p22.Range.End)
$afterPrint.Find.Execute("x_center", $false, $false, $false, $false, $false, $true, 1, $false, "x_center", 2)

# ---------------------------------------------------------------------
# 6) ") (line1" + "23" + "," + "124" + ")" -> ") (line123,124)"
# ---------------------------------------------------------------------
$p22.Range.Find.Execute(") (line123,124)", $false, $false, $false, $false, $false, $true, 1, $false, ") (line123,124)", 2)

# ---------------------------------------------------------------------
# 7) "print(type(distance)) print(distance) Provides distance of the
#    object from the camera" + "." -> single merged run
# ---------------------------------------------------------------------
$p39 = $d.Paragraphs.Item(39)
$p39.Range.Find.Execute("print(type(distance)) print(distance) Provides distance of the object from the camera.", $false, $false, $false, $false, $false, $true, 1, $false, "print(type(distance)) print(distance) Provides distance of the object from the camera.", 2)

# ---------------------------------------------------------------------
# 8) Insert two new paragraphs right after the "5.58368138988692"
#    paragraph (before the existing blank ListParagraph that follows
#    it): one blank ListParagraph, and one ListParagraph with the new
#    explanatory text about how distance is calculated.
# ---------------------------------------------------------------------
$p45 = $d.Paragraphs.Item(45)
$insertPoint = $p45.Range
$insertPoint.Collapse(0)
$insertPoint.InsertParagraphAfter()

$pBlank = $d.Paragraphs.Item(46)
$pBlank.Style = "ListParagraph"
$blankEnd = $pBlank.Range
$blankEnd.Collapse(0)
$blankEnd.InsertParagraphAfter()

$pText = $d.Paragraphs.Item(47)
$pText.Style = "ListParagraph"
$pText.Range.Text = "The distance is calculated as sqroot of x2+y2+z2. This is because x,y,z represent the coordinate values and the "

# Mark "sqroot" and "x,y" / ",z" as spell/gram checked the same way the
# original document does elsewhere (adds the proofErr markers around
# them to match the authored text precisely).
$sqrootRng = $pText.Range
$sqrootRng.Find.Execute("sqroot", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# ---------------------------------------------------------------------
# 9) " = " + "5" -> " = 5"
# ---------------------------------------------------------------------
$p55 = $d.Paragraphs.Item($d.Paragraphs.Count)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -like "At research_distance_box = 5*") {
        $pp.Range.Find.Execute("At research_distance_box = 5", $false, $false, $false, $false, $false, $true, 1, $false, "At research_distance_box = 5", 2)
        break
    }
}

Write-Output "done"
